# Scenario_EnergyCarrier_Price_MarkUp.xlsx - "cost precalculated in scenario"
# Extend the id_energy_carrier (column D) list from [1,2,6,8,12,13,14,25] to the
# full [1,2,3,6,8,12,13,14,15,19,25,26] set for both id_sector groups (3 and 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target list of id_energy_carrier values, per id_sector group.
$dVals = @(1,2,3,6,8,12,13,14,15,19,25,26)

# --- Rows 2-13: id_sector = 3 -------------------------------------------
$startRow = 2
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = 3
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}

# --- Rows 14-25: id_sector = 6 -------------------------------------------
# Rows 14-17 already exist (need their C/D values updated); rows 18-25 are
# brand new and must be created (values + number-format copied from row 17).
$startRow2 = 14
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $r = $startRow2 + $i

    if ($r -gt 17) {
        # New row: clone formatting (incl. the 0.000 number format on F:AK)
        # from the last existing data row (17) before filling in values.
        $ws.Range("A17:AK17").Copy()
        $ws.Range("A" + $r + ":AK" + $r).PasteSpecial(-4122)

        $ws.Cells.Item($r, 1).Value = 1
        $ws.Cells.Item($r, 2).Value = 9
        $ws.Cells.Item($r, 5).Value = "euro/kWh"
        for ($col = 6; $col -le 37; $col++) {
            $ws.Cells.Item($r, $col).Value = 0
        }
    }

    $ws.Cells.Item($r, 3).Value = 6
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}

$excel.CutCopyMode = 0

# Resize the table/ListObject + autofilter to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:AK25"))

# Reposition the selection / view the way the saved file shows it.
$ws.Cells.Item(18, 5).Select()
$excel.ActiveWindow.Zoom = 88
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
